$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("dSF") updates per repulled/pushed data and mean recalculation
$updates = @{
    8  = -8
    10 = -6
    14 = -2
    19 = -4
    20 = -10
    25 = -5
    26 = -2
    27 = -8
    28 = 4
    31 = -1
    34 = -2
    36 = 9
    37 = 1
    39 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
